# Apply corrections to event names / capitalization and fill in the
# missing "number of people per team" value for the last event row so the
# team-count algorithm can pick it up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix a mis-typed event name in the events list (column F)
$ws.Range("F7").Value = "Optics"

# Normalize capitalization / casing of event names inside the
# "Events I want to do:" column (column C)
$ws.Range("C8").Value = "Fossils, Material Science, Robot Tour"
$ws.Range("C11").Value = "Experimental Design, Optics"
$ws.Range("C12").Value = "Optics, Write It Do It"
$ws.Range("C15").Value = "Write It Do It, Electric Vehicle"
$ws.Range("C16").Value = "Electric Vehicle, Dynamic Planet"

# The last event (Microbe Mission, row 24) was missing its
# "number of people per team" value -- fill it in like the other events.
$ws.Range("J23").Copy($ws.Range("J24"))
$ws.Range("J24").Value = 2.0
